# Edit script: add ValidationData hidden sheet, restructure data validations
# to reference it, and bold the header row on the main sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the hidden "ValidationData" worksheet right after the main sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ValidationData"

# --- Populate the lookup columns (row-major fill order matches the source data) ---
$ws2.Range("E1").Value = 'form'
$ws2.Range("F1").Value = 'Pre-seed'
$ws2.Range("G1").Value = "-- Tech`n    B2B SaaS"
$ws2.Range("H1").Value = 'Global'
$ws2.Range("I1").Value = 'simple'
$ws2.Range("J1").Value = '1-5'
$ws2.Range("K1").Value = 'pitch_deck'
$ws2.Range("N1").Value = 'FREE'
$ws2.Range("E2").Value = 'email'
$ws2.Range("F2").Value = 'Seed'
$ws2.Range("G2").Value = 'Fintech'
$ws2.Range("H2").Value = 'North America'
$ws2.Range("I2").Value = 'standard'
$ws2.Range("J2").Value = '6-10'
$ws2.Range("K2").Value = 'video'
$ws2.Range("N2").Value = 'PRO'
$ws2.Range("E3").Value = 'other'
$ws2.Range("F3").Value = 'Series A'
$ws2.Range("G3").Value = 'Healthtech'
$ws2.Range("H3").Value = 'South America'
$ws2.Range("I3").Value = 'comprehensive'
$ws2.Range("J3").Value = '11-20'
$ws2.Range("K3").Value = 'financial_projections'
$ws2.Range("N3").Value = 'MAX'
$ws2.Range("F4").Value = 'Series B'
$ws2.Range("G4").Value = 'AI/ML'
$ws2.Range("H4").Value = 'LATAM'
$ws2.Range("J4").Value = '21+'
$ws2.Range("K4").Value = 'business_plan'
$ws2.Range("F5").Value = 'Series C'
$ws2.Range("G5").Value = 'Deep tech'
$ws2.Range("H5").Value = 'Europe'
$ws2.Range("K5").Value = 'traction_data'
$ws2.Range("F6").Value = 'Growth'
$ws2.Range("G6").Value = 'Climate tech'
$ws2.Range("H6").Value = 'Western Europe'
$ws2.Range("F7").Value = 'All'
$ws2.Range("G7").Value = 'Consumer'
$ws2.Range("H7").Value = 'Eastern Europe'
$ws2.Range("G8").Value = 'E-commerce'
$ws2.Range("H8").Value = 'Continental Europe'
$ws2.Range("G9").Value = 'Marketplace'
$ws2.Range("H9").Value = 'Middle East'
$ws2.Range("G10").Value = 'Gaming'
$ws2.Range("H10").Value = 'Africa'
$ws2.Range("G11").Value = 'Web3'
$ws2.Range("H11").Value = 'Asia'
$ws2.Range("G12").Value = 'Developer tools'
$ws2.Range("H12").Value = 'East Asia'
$ws2.Range("G13").Value = 'Cybersecurity'
$ws2.Range("H13").Value = 'South Asia'
$ws2.Range("G14").Value = 'Logistics'
$ws2.Range("H14").Value = 'South East Asia'
$ws2.Range("G15").Value = 'AdTech'
$ws2.Range("H15").Value = 'Oceania'
$ws2.Range("G16").Value = 'PropTech'
$ws2.Range("H16").Value = 'EMEA'
$ws2.Range("G17").Value = 'InsurTech'
$ws2.Range("H17").Value = 'Emerging Markets'
$ws2.Range("G18").Value = "-- Non-Tech / Other`n    Agriculture"
$ws2.Range("G19").Value = 'Automotive'
$ws2.Range("G20").Value = 'Biotechnology'
$ws2.Range("G21").Value = 'Construction'
$ws2.Range("G22").Value = 'Consulting'
$ws2.Range("G23").Value = 'Consumer Goods'
$ws2.Range("G24").Value = 'Education'
$ws2.Range("G25").Value = 'Energy'
$ws2.Range("G26").Value = 'Entertainment'
$ws2.Range("G27").Value = 'Environmental Services'
$ws2.Range("G28").Value = 'Fashion'
$ws2.Range("G29").Value = 'Food & Beverage'
$ws2.Range("G30").Value = 'Government'
$ws2.Range("G31").Value = 'Healthcare Services'
$ws2.Range("G32").Value = 'Hospitality'
$ws2.Range("G33").Value = 'Human Resources'
$ws2.Range("G34").Value = 'Insurance'
$ws2.Range("G35").Value = 'Legal'
$ws2.Range("G36").Value = 'Manufacturing'
$ws2.Range("G37").Value = 'Media'
$ws2.Range("G38").Value = 'Non-profit'
$ws2.Range("G39").Value = 'Pharmaceuticals'
$ws2.Range("G40").Value = 'Real Estate'
$ws2.Range("G41").Value = 'Retail'
$ws2.Range("G42").Value = 'Telecommunications'
$ws2.Range("G43").Value = 'Transportation'
$ws2.Range("G44").Value = 'Utilities'
$ws2.Range("G45").Value = 'Other'

$ws2.Visible = $false

# --- Bold the header row on the main sheet ---
$ws1.Range("A1:N1").Font.Bold = $true

# --- Re-point each data validation rule at the ValidationData sheet ranges,
#     and refresh the error message text to match.
#     (Operator args are passed as [Missing] so the serialized XML omits
#     an explicit operator="equal" attribute, matching a plain list rule.) ---

$rng = $ws1.Range("E10:E1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$E`$1:`$E`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("E2:E1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$E`$1:`$E`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("F10:F1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$F`$1:`$F`$7")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("F2:F1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$F`$1:`$F`$7")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("G10:G1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$G`$1:`$G`$45")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("G2:G1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$G`$1:`$G`$45")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("H10:H1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$H`$1:`$H`$17")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("H2:H1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$H`$1:`$H`$17")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("I10:I1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$I`$1:`$I`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("I2:I1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$I`$1:`$I`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("J10:J1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$J`$1:`$J`$4")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("J2:J1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$J`$1:`$J`$4")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("K10:K1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$K`$1:`$K`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("K2:K1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$K`$1:`$K`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("N10:N1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$N`$1:`$N`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

$rng = $ws1.Range("N2:N1000")
$v = $rng.Validation
$v.Modify(3, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "=ValidationData!`$N`$1:`$N`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

# --- Leave the main sheet as the active/selected one ---
$ws1.Activate()

Write-Output "done"
